# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets to reflect a refreshed data scrape.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1295
$ws.Range("F5").Value = 7461
$ws.Range("F6").Value = 1867
$ws.Range("F7").Value = 6445
$ws.Range("F8").Value = 150
$ws.Range("F9").Value = 2009
$ws.Range("F10").Value = 541
$ws.Range("F11").Value = 33
$ws.Range("F13").Value = 42
$ws.Range("F16").Value = 59
$ws.Range("F17").Value = 8297
$ws.Range("F18").Value = 150
$ws.Range("F19").Value = 60
$ws.Range("F20").Value = 192
$ws.Range("F22").Value = 1786
$ws.Range("F30").Value = 1947
$ws.Range("F31").Value = 834
$ws.Range("F32").Value = 445
$ws.Range("F33").Value = 5
$ws.Range("F35").Value = 148
$ws.Range("F37").Value = 95
$ws.Range("F38").Value = 3945

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 384
$ws.Range("F12").Value = 10

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2307
$ws.Range("F4").Value = 699
$ws.Range("F5").Value = 296

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2307
$ws.Range("F4").Value = 1295
$ws.Range("F5").Value = 384
$ws.Range("F6").Value = 7461
$ws.Range("F8").Value = 1867
$ws.Range("F9").Value = 6445
$ws.Range("F10").Value = 2009
$ws.Range("F13").Value = 541
$ws.Range("F14").Value = 33
$ws.Range("F18").Value = 42
$ws.Range("F20").Value = 10
$ws.Range("F22").Value = 59
$ws.Range("F23").Value = 8297
$ws.Range("F24").Value = 60
$ws.Range("F25").Value = 192
$ws.Range("F27").Value = 1786
$ws.Range("F32").Value = 1947
$ws.Range("F33").Value = 834
$ws.Range("F35").Value = 445
$ws.Range("F36").Value = 5
$ws.Range("F40").Value = 148
$ws.Range("F42").Value = 95
